# EPBDS-9540 Support Java Name convestion on Json field name generating in
# SpreadsheetResults. Rework.
#
# The "Rules" sheet has a block of expectation-formula-looking text cells
# (C50:C65) that reference Step-map keys using the *lowercase* key names
# (e.g. "step2", "step4", ...). Rework them to use the *capitalized* Java
# field-name convention ("Step2", "Step4", ...), matching the keys now
# produced by the generator. The cells hold literal text beginning with
# "=" (quote-prefixed), not real formulas, so each new value is entered
# with a leading apostrophe to keep it literal text instead of having
# Excel parse it as a formula.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

$ws.Range("C50").Value = "'= instanceOf(`$Step2[""Step2""], Map.class)"
$ws.Range("C51").Value = "'= instanceOf(`$Step2[""Step4""], java.lang.reflect.Array.newInstance(Map.class, 1).getClass())"
$ws.Range("C52").Value = "'= instanceOf(((Object[])`$Step2[""Step8""])[0], Map.class)"
$ws.Range("C53").Value = "'=  instanceOf(((Object[][])`$Step2[""Step6""])[0][0], Map.class)"
$ws.Range("C54").Value = "'= instanceOf(((Object[])`$Step2[""Step14""])[0], Map.class)"
$ws.Range("C55").Value = "'= instanceOf(((Object[][])`$Step2[""Step15""])[0][0], Map.class)"
$ws.Range("C56").Value = "'= instanceOf(`$Step2[""Step13""], Map.class)"
$ws.Range("C57").Value = "'= instanceOf(((Object[][])`$Step2[""Step9""])[0][0], Map.class)"
$ws.Range("C58").Value = "'=  instanceOf(((Object[][])`$Step2[""Step10""])[0][0], Map.class)"
$ws.Range("C59").Value = "'= instanceOf(((Object[])`$Step2[""Step11""])[0], Map.class)"
$ws.Range("C60").Value = "'= instanceOf(((Map)`$Step2[""Step16""]).keySet().iterator().next(), Map.class)"
$ws.Range("C61").Value = "'= instanceOf(((Map)`$Step2[""Step16""]).values().iterator().next(), Map.class)"
$ws.Range("C62").Value = "'= instanceOf(((List)`$Step2[""Step17""])[0], Map.class)"
$ws.Range("C63").Value = "'= instanceOf(((Collection)`$Step2[""Step18""]).iterator().next(), Map.class)"
$ws.Range("C64").Value = "'= instanceOf(((Collection)`$Step2[""Step19""]).iterator().next(), Map.class)"
$ws.Range("C65").Value = "'= instanceOf(`$Step2[""Step20""], Map.class)"

# Reflects the author's last cursor position on the "Rules" sheet when the
# edits were saved.
$ws.Activate()
$ws.Range("J57").Select()
